$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (october-2025)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "november-2025"

# Populate A1 with the new "Total Gross Cumulative Voted Spending" line for November 2025
$newSheet.Range("A1").Value = " Total Gross Cumulative Voted Spending                               96,914       97,288          374          0.4%               5,240         5.7%             107,015       103,472        3,543         3.4%"

# Match the print/page setup used by the other monthly sheets
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1
